# Incorporating bat exogenous and endogenous taxa
# Adds a new data row (row 7) to the "side_data.txt" worksheet for the
# Eptesicus fuscus deltaretrovirus (EtDRV) entry, mirroring the layout
# and formatting of the existing taxon rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("side_data.txt")
$ws.Activate() | Out-Null

# --- Copy row 6's formatting down into the new row 7 -----------------------
$ws.Range("A6:I6").Copy() | Out-Null
$ws.Range("A7:I7").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# Columns B and C on the new row use the plain data style (style index 1,
# matching rows 2-5) rather than the "last row" accent style (index 4) that
# row 6 uses for those columns - copy that formatting across too.
$ws.Range("B2:C2").Copy() | Out-Null
$ws.Range("B7:C7").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# --- Populate the new row's values ------------------------------------------
$ws.Range("A7").Value = "MT700539"
$ws.Range("B7").Value = "EtDRV"
$ws.Range("C7").Value = "Eptesicus fuscus deltaretrovirus"
$ws.Range("D7").Value = "Orthoretrovirinae"
$ws.Range("E7").Value = "Clade II"
$ws.Range("F7").Value = "Deltaretrovirus"
$ws.Range("G7").Value = "Primate"
$ws.Range("H7").Value = "Eptesicus fuscus"
$ws.Range("I7").Value = "Big brown bat"

# --- Update the view selection to cover the newly extended table -----------
$ws.Range("A1:I7").Select() | Out-Null

Write-Host "Added EtDRV (MT700539) row to side_data.txt"
